$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> [new nombre_aides (col C), new montant_total (col E)]
$updates = @{
    2   = @(766351, 1429263437)
    13  = @(187874, 1168847701)
    19  = @(27520, 132568759)
    21  = @(175244, 316831538)
    41  = @(126947, 662714913)
    57  = @(31599, 162636444)
    81  = @(88360, 499741973)
    88  = @(71284, 110329804)
    121 = @(1306505, 2275852587)
    129 = @(633969, 3437685258)
    132 = @(586157, 3476085690)
    151 = @(39937, 60395303)
    156 = @(12418, 40886248)
    171 = @(95834, 490727046)
    178 = @(515894, 891232265)
    186 = @(236848, 1190253497)
    237 = @(283332, 1438703447)
    240 = @(205945, 1070200493)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
